# "Generate Report for Handback": the handback XLIFFs have come back in
# sync with en-US, so the status/report workbook is updated to reflect it
# for both language sheets (zh-cn, de-de), and the two report tables grow
# three new "populated" columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) worth of real width.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdHyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0ec6b17b9b3728738dd9fbd6571fc1cbae04bf4a/e2e/4a63b18c-daab-41cd-84bc-2f45df94d6b0.md"

# --- Overview sheet: widen the two per-language status columns (E, F),
# and roll the same "handed back" status up into its summary cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# --- Per-language report sheets: zh-cn and de-de got handed back ---
$langSheets = @("zh-cn", "de-de")

foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)

    # Status column (C) got wider to fit the longer status text, and the
    # newly-populated Latest Target File / Latest Handback File columns
    # (I, J) now need full 40-char width like the other file-name columns.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Status: handoff -> handback, now in sync with en-US.
    $ws.Range("C2").Value = $statusText

    # Latest Target File (I2): the localized file is now the same source
    # markdown document referenced by A2, with the same kind of hyperlink.
    $sourceName = $ws.Range("A2").Value()
    $ws.Range("I2").Value = $sourceName
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdHyperlinkTarget, "", "", $sourceName) | Out-Null

    # Latest Handback File (J2): the xliff that was handed off is the one
    # that was handed back (G2 already holds that filename).
    $ws.Range("J2").Value = $ws.Range("G2").Value()
}

# Latest Handback DateTime (K2): each language got its own handback stamp.
$wb.Worksheets.Item("zh-cn").Range("K2").Value = "2016-08-17 12:56:33"
$wb.Worksheets.Item("de-de").Range("K2").Value = "2016-08-17 12:56:41"
